$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 23810010
$ws.Range("J19").Value = 35714884
$ws.Range("L19").Value = 35714884
$ws.Range("N19").Value = -35715234
$ws.Range("H69").Value = 5234.2856
$ws.Range("I69").Value = 5792.143
$ws.Range("K69").Value = 17376.429
$ws.Range("M69").Value = -16502.429
$ws.Range("H72").Value = 5234.2856
$ws.Range("I72").Value = 5792.143
$ws.Range("K72").Value = 52129.287
$ws.Range("M72").Value = -47761.287
$ws.Range("H107").Value = 896
$ws.Range("I107").Value = 820.4545000000001
$ws.Range("J107").Value = 999.875
$ws.Range("K107").Value = 820.4545000000001
$ws.Range("L107").Value = 999.875
$ws.Range("M107").Value = 1099.5455
$ws.Range("N107").Value = -4839.875
$ws.Range("H129").Value = 1542.9688
$ws.Range("I129").Value = 868.5714
$ws.Range("J129").Value = 1731.8
$ws.Range("K129").Value = 2605.7142
$ws.Range("L129").Value = 5195.4
$ws.Range("M129").Value = 2394.2858
$ws.Range("N129").Value = -15195.4
$ws.Range("H134").Value = 111184180
$ws.Range("J134").Value = 111184180
$ws.Range("L134").Value = 111184180
$ws.Range("N134").Value = -111194320

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 22588.4
$ws.Range("I31").Value = 22588.4
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 22588.4
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -22294.4
$ws.Range("N31").ClearContents()
$ws.Range("H32").Value = 7955.24
$ws.Range("I32").Value = 6030.6
$ws.Range("J32").Value = 25277
$ws.Range("K32").Value = 6030.6
$ws.Range("L32").Value = 25277
$ws.Range("M32").Value = -5743.6
$ws.Range("N32").Value = -25851
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H102").Value = 2865.5557
$ws.Range("I102").Value = 2335
$ws.Range("J102").Value = 3926.6667
$ws.Range("K102").Value = 2335
$ws.Range("L102").Value = 3926.6667
$ws.Range("M102").Value = -713
$ws.Range("N102").Value = -7170.6667
$ws.Range("H110").Value = 1546.7646
$ws.Range("I110").Value = 1630.2
$ws.Range("J110").Value = 1427.5714
$ws.Range("K110").Value = 1630.2
$ws.Range("L110").Value = 1427.5714
$ws.Range("M110").Value = 414.8
$ws.Range("N110").Value = -5517.5714
$ws.Range("H122").Value = 1168.1052
$ws.Range("I122").Value = 1022.61536
$ws.Range("J122").Value = 1483.3334
$ws.Range("K122").Value = 3067.84608
$ws.Range("L122").Value = 4450.0002
$ws.Range("M122").Value = -617.8460800000003
$ws.Range("N122").Value = -9350.0002
$ws.Range("H132").Value = 1259106
$ws.Range("I132").Value = 2074083.2
$ws.Range("J132").Value = 9474.134
$ws.Range("K132").Value = 6222249.6
$ws.Range("L132").Value = 28422.402
$ws.Range("M132").Value = -6219719.6
$ws.Range("N132").Value = -33482.402

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 955.2308
$ws.Range("I94").Value = 684.1739
$ws.Range("J94").Value = 3033.3333
$ws.Range("K94").Value = 684.1739
$ws.Range("L94").Value = 3033.3333
$ws.Range("M94").Value = -233.1739
$ws.Range("N94").Value = -3935.3333
$ws.Range("H98").Value = 25000
$ws.Range("J98").Value = 25000
$ws.Range("L98").Value = 25000
$ws.Range("N98").Value = -30990
$ws.Range("H134").Value = 4956.7446
$ws.Range("I134").Value = 1720.2273
$ws.Range("J134").Value = 7804.88
$ws.Range("K134").Value = 5160.6819
$ws.Range("L134").Value = 23414.64
$ws.Range("M134").Value = -2625.6819
$ws.Range("N134").Value = -28484.64

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1224085.8
$ws.Range("I58").Value = 2971.75
$ws.Range("J58").Value = 2948011.5
$ws.Range("K58").Value = 2971.75
$ws.Range("L58").Value = 2948011.5
$ws.Range("M58").Value = -2768.75
$ws.Range("N58").Value = -2948417.5
$ws.Range("H132").Value = 2639.8918
$ws.Range("I132").Value = 1931.9166
$ws.Range("J132").Value = 3946.923
$ws.Range("K132").Value = 5795.7498
$ws.Range("L132").Value = 11840.769
$ws.Range("M132").Value = -3265.7498
$ws.Range("N132").Value = -16900.769
$ws.Range("H134").Value = 2002.56
$ws.Range("I134").Value = 1208.3636
$ws.Range("J134").Value = 7826.6665
$ws.Range("K134").Value = 3625.0908
$ws.Range("L134").Value = 23479.9995
$ws.Range("M134").Value = -1090.0908
$ws.Range("N134").Value = -28549.9995
$ws.Range("H136").Value = 1224085.8
$ws.Range("I136").Value = 2971.75
$ws.Range("J136").Value = 2948011.5
$ws.Range("K136").Value = 8915.25
$ws.Range("L136").Value = 8844034.5
$ws.Range("M136").Value = -6365.25
$ws.Range("N136").Value = -8849134.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 8649.833000000001
$ws.Range("I4").Value = 12724.75
$ws.Range("J4").Value = 500
$ws.Range("K4").Value = 38174.25
$ws.Range("L4").Value = 1500
$ws.Range("M4").Value = -38062.25
$ws.Range("N4").Value = -1724
$ws.Range("H34").Value = 1463.6364
$ws.Range("J34").Value = 1700
$ws.Range("L34").Value = 5100
$ws.Range("N34").Value = -5268
$ws.Range("H50").Value = 116.625
$ws.Range("I50").Value = 61.857143
$ws.Range("K50").Value = 185.571429
$ws.Range("M50").Value = 295.428571
$ws.Range("H53").Value = 116.625
$ws.Range("I53").Value = 61.857143
$ws.Range("K53").Value = 185.571429
$ws.Range("M53").Value = 295.428571
$ws.Range("H131").Value = 1080.9149
$ws.Range("I131").Value = 583.3333
$ws.Range("J131").Value = 1153.7317
$ws.Range("K131").Value = 1749.9999
$ws.Range("L131").Value = 3461.1951
$ws.Range("M131").Value = 3290.0001
$ws.Range("N131").Value = -13541.1951
$ws.Range("H137").Value = 2621.8845
$ws.Range("I137").Value = 1299.375
$ws.Range("J137").Value = 3209.6667
$ws.Range("K137").Value = 3898.125
$ws.Range("L137").Value = 9629.000100000001
$ws.Range("M137").Value = 1201.875
$ws.Range("N137").Value = -19829.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 10840.363
$ws.Range("J95").Value = 10840.363
$ws.Range("L95").Value = 10840.363
$ws.Range("N95").Value = -16332.363
$ws.Range("H132").Value = 2674.75
$ws.Range("I132").Value = 1734.6428
$ws.Range("K132").Value = 5203.928400000001
$ws.Range("M132").Value = -2673.928400000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 355.6
$ws.Range("I55").Value = 359.8889
$ws.Range("J55").Value = 349.16666
$ws.Range("K55").Value = 359.8889
$ws.Range("L55").Value = 349.16666
$ws.Range("M55").Value = -186.8889
$ws.Range("N55").Value = -695.16666
$ws.Range("H97").Value = 12836.462
$ws.Range("J97").Value = 12836.462
$ws.Range("L97").Value = 12836.462
$ws.Range("N97").Value = -14818.462

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 70991.60000000001
$ws.Range("J46").Value = 70991.60000000001
$ws.Range("L46").Value = 70991.60000000001
$ws.Range("N46").Value = -71453.60000000001
$ws.Range("H81").Value = 2500.2727
$ws.Range("I81").Value = 2550.3
$ws.Range("J81").Value = 2000
$ws.Range("K81").Value = 5100.6
$ws.Range("L81").Value = 4000
$ws.Range("M81").Value = -4039.6
$ws.Range("N81").Value = -6122
$ws.Range("H84").Value = 2500.2727
$ws.Range("I84").Value = 2550.3
$ws.Range("J84").Value = 2000
$ws.Range("K84").Value = 25503
$ws.Range("L84").Value = 20000
$ws.Range("M84").Value = -20199
$ws.Range("N84").Value = -30608
$ws.Range("H97").Value = 20078
$ws.Range("J97").Value = 20078
$ws.Range("L97").Value = 20078
$ws.Range("N97").Value = -22060
$ws.Range("H134").Value = 70991.60000000001
$ws.Range("J134").Value = 70991.60000000001
$ws.Range("L134").Value = 212974.8
$ws.Range("N134").Value = -218044.8
